$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp in cell A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 16:13"

# Update COVID data values (columns B-H) for the rows whose figures changed
$ws.Range("B4").Value = 7726175
$ws.Range("C4").Value = 3429
$ws.Range("E4").Value = 2573411
$ws.Range("G4").Value = 87
$ws.Range("H4").Value = 215909
$ws.Range("B5").Value = 6764710
$ws.Range("C5").Value = 10531
$ws.Range("D5").Value = 5750403
$ws.Range("E5").Value = 909656
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 104651
$ws.Range("B18").Value = 391044
$ws.Range("C18").Value = 3923
$ws.Range("D18").Value = 319784
$ws.Range("E18").Value = 61656
$ws.Range("G18").Value = 73
$ws.Range("H18").Value = 9604
$ws.Range("B26").Value = 308340
$ws.Range("C26").Value = 1221
$ws.Range("E26").Value = 33098
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = 9642
$ws.Range("B53").Value = 81256
$ws.Range("C53").Value = 944
$ws.Range("D53").Value = 51037
$ws.Range("E53").Value = 28179
$ws.Range("G53").Value = 8
$ws.Range("H53").Value = 2040
$ws.Range("B54").Value = 80662
$ws.Range("C54").Value = 642
$ws.Range("D54").Value = 30131
$ws.Range("E54").Value = 48084
$ws.Range("G54").Value = 14
$ws.Range("H54").Value = 2447
$ws.Range("B70").Value = 42840
$ws.Range("C70").Value = 408
$ws.Range("D70").Value = 35953
$ws.Range("E70").Value = 6532
$ws.Range("G70").Value = 6
$ws.Range("H70").Value = 355
$ws.Range("B76").Value = 34193
$ws.Range("C76").Value = 121
$ws.Range("E76").Value = 1899
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 758
$ws.Range("B79").Value = 29737
$ws.Range("C79").Value = 103
$ws.Range("D79").Value = 24643
$ws.Range("E79").Value = 4221
$ws.Range("B80").Value = 29075
$ws.Range("C80").Value = 365
$ws.Range("D80").Value = 22614
$ws.Range("E80").Value = 5553
$ws.Range("G80").Value = 5
$ws.Range("H80").Value = 908
$ws.Range("B85").Value = 21433
$ws.Range("C85").Value = 1400
$ws.Range("D85").Value = 6084
$ws.Range("E85").Value = 14839
$ws.Range("G85").Value = 39
$ws.Range("H85").Value = 510
$ws.Range("B86").Value = 20924
$ws.Range("D86").Value = 19764
$ws.Range("E86").Value = 740
$ws.Range("B87").Value = 20541
$ws.Range("D87").Value = 9989
$ws.Range("E87").Value = 10132
$ws.Range("H87").Value = 420
$ws.Range("B89").Value = 19413
$ws.Range("C89").Value = 317
$ws.Range("D89").Value = 15749
$ws.Range("E89").Value = 2892
$ws.Range("G89").Value = 4
$ws.Range("H89").Value = 772
$ws.Range("B95").Value = 14895
$ws.Range("C95").Value = 111
$ws.Range("D95").Value = 11863
$ws.Range("E95").Value = 2757
$ws.Range("B107").Value = 10055
$ws.Range("C107").Value = 41
$ws.Range("D107").Value = 8876
$ws.Range("E107").Value = 1101
$ws.Range("B148").Value = 3172
$ws.Range("C148").Value = 91
$ws.Range("D148").Value = 2366
$ws.Range("E148").Value = 796
$ws.Range("B196").Value = 131
$ws.Range("C196").Value = 1
$ws.Range("E196").Value = 14
$ws.Range("B214").Value = 15
$ws.Range("C214").Value = 1
$ws.Range("E214").Value = 1
